$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the old row 625 (the "2026/12/29" block),
# shifting every subsequent row down by 2 (625->627 ... 666->668).
$ws.Rows("625:626").Insert()

# New row 625: continuation of the 2026/01/14 (Wed) block, hour 22
$ws.Range("A625").NumberFormat = "@"
$ws.Range("A625").Value = "2026/01/14"
$ws.Range("A625").Style = "Normal"
$ws.Range("B625").Value = "水"
$ws.Range("C625").Value = 22
$ws.Range("D625").Value = 28

# New row 626: first entry of 2026/01/15 (Thu), hour 2
$ws.Range("A626").NumberFormat = "@"
$ws.Range("A626").Value = "2026/01/15"
$ws.Range("A626").Style = "Normal"
$ws.Range("B626").Value = "木"
$ws.Range("C626").Value = 2
$ws.Range("D626").Value = 29
